$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 4 cell(s) changed
$ws.Range("B2").Value2 = 'd83b6dd9-45c7-4ce5-b299-6086a4e5ba79'
$ws.Range("C2").Value2 = '2024-07-25T05:44:00.000Z'
$ws.Range("D2").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I2").Value2 = 'https://www.notion.so/1-d83b6dd945c74ce5b2996086a4e5ba79'

# Row 3: 4 cell(s) changed
$ws.Range("B3").Value2 = '77981efc-5085-4ceb-8daf-433ff69a8790'
$ws.Range("C3").Value2 = '2024-07-23T09:34:00.000Z'
$ws.Range("D3").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I3").Value2 = 'https://www.notion.so/2-77981efc50854ceb8daf433ff69a8790'

# Row 4: 4 cell(s) changed
$ws.Range("B4").Value2 = '25f7c61d-806e-44e4-b06a-d1cdd69e0209'
$ws.Range("C4").Value2 = '2024-07-21T10:08:00.000Z'
$ws.Range("D4").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I4").Value2 = 'https://www.notion.so/3-25f7c61d806e44e4b06ad1cdd69e0209'

# Row 5: 4 cell(s) changed
$ws.Range("B5").Value2 = '6e31e972-3b9d-4943-9bca-9623487e0a7d'
$ws.Range("C5").Value2 = '2024-07-20T09:48:00.000Z'
$ws.Range("D5").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I5").Value2 = 'https://www.notion.so/4-6e31e9723b9d49439bca9623487e0a7d'

# Row 6: 4 cell(s) changed
$ws.Range("B6").Value2 = '94fa7304-3963-4bcd-8356-5abbfce2a6f3'
$ws.Range("C6").Value2 = '2024-07-19T10:42:00.000Z'
$ws.Range("D6").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I6").Value2 = 'https://www.notion.so/5-94fa730439634bcd83565abbfce2a6f3'

# Row 7: 4 cell(s) changed
$ws.Range("B7").Value2 = '8a14bd80-4464-4e14-9c3b-b8b92e5e34f1'
$ws.Range("C7").Value2 = '2024-07-18T10:59:00.000Z'
$ws.Range("D7").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I7").Value2 = 'https://www.notion.so/6-8a14bd8044644e149c3bb8b92e5e34f1'

# Row 8: 3 cell(s) changed
$ws.Range("B8").Value2 = 'ab7f95ce-43b3-42b1-83b5-0578395cef2d'
$ws.Range("D8").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I8").Value2 = 'https://www.notion.so/7-ab7f95ce43b342b183b50578395cef2d'

# Row 9: 4 cell(s) changed
$ws.Range("B9").Value2 = '3b58dea8-d569-439c-9d28-37bc5cb00b87'
$ws.Range("C9").Value2 = '2024-07-17T12:17:00.000Z'
$ws.Range("D9").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I9").Value2 = 'https://www.notion.so/8-3b58dea8d569439c9d2837bc5cb00b87'

# Row 10: 4 cell(s) changed
$ws.Range("B10").Value2 = '0fecdb7c-3c37-4eea-9728-9e9e308f8f91'
$ws.Range("C10").Value2 = '2024-07-15T04:38:00.000Z'
$ws.Range("D10").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I10").Value2 = 'https://www.notion.so/9-0fecdb7c3c374eea97289e9e308f8f91'

# Row 11: 4 cell(s) changed
$ws.Range("B11").Value2 = 'e64e264f-45f0-48db-8349-7eb369265ce0'
$ws.Range("C11").Value2 = '2024-07-14T14:34:00.000Z'
$ws.Range("D11").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I11").Value2 = 'https://www.notion.so/10-e64e264f45f048db83497eb369265ce0'

# Row 12: 4 cell(s) changed
$ws.Range("B12").Value2 = '6c9a4fbe-94ba-4813-a444-2725e44cd033'
$ws.Range("C12").Value2 = '2024-07-13T09:28:00.000Z'
$ws.Range("D12").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I12").Value2 = 'https://www.notion.so/11-6c9a4fbe94ba4813a4442725e44cd033'

# Row 13: 4 cell(s) changed
$ws.Range("B13").Value2 = '5cd2df68-475f-4500-b2d7-678411a240a9'
$ws.Range("C13").Value2 = '2024-07-12T09:53:00.000Z'
$ws.Range("D13").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I13").Value2 = 'https://www.notion.so/12-5cd2df68475f4500b2d7678411a240a9'

# Row 14: 4 cell(s) changed
$ws.Range("B14").Value2 = '9afdbbb6-3c71-4c0c-9dde-e729964648b8'
$ws.Range("C14").Value2 = '2024-07-11T06:12:00.000Z'
$ws.Range("D14").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I14").Value2 = 'https://www.notion.so/13-9afdbbb63c714c0c9ddee729964648b8'

# Row 15: 4 cell(s) changed
$ws.Range("B15").Value2 = 'fb2d1794-0cdc-46f0-a9a1-07162fa1770c'
$ws.Range("C15").Value2 = '2024-07-09T08:22:00.000Z'
$ws.Range("D15").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I15").Value2 = 'https://www.notion.so/14-fb2d17940cdc46f0a9a107162fa1770c'

# Row 16: 4 cell(s) changed
$ws.Range("B16").Value2 = '107624f7-7695-4a1b-b0ba-d8f486d99ffc'
$ws.Range("C16").Value2 = '2024-07-08T04:11:00.000Z'
$ws.Range("D16").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I16").Value2 = 'https://www.notion.so/15-107624f776954a1bb0bad8f486d99ffc'

# Row 17: 4 cell(s) changed
$ws.Range("B17").Value2 = '77beaf39-139f-468c-8b7a-9dfbffdc34ce'
$ws.Range("C17").Value2 = '2024-07-06T08:51:00.000Z'
$ws.Range("D17").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I17").Value2 = 'https://www.notion.so/16-77beaf39139f468c8b7a9dfbffdc34ce'

# Row 18: 4 cell(s) changed
$ws.Range("B18").Value2 = '03447847-60f7-4bcc-9b55-b672299a1911'
$ws.Range("C18").Value2 = '2024-07-04T03:57:00.000Z'
$ws.Range("D18").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I18").Value2 = 'https://www.notion.so/17-0344784760f74bcc9b55b672299a1911'

# Row 19: 4 cell(s) changed
$ws.Range("B19").Value2 = 'd37aca60-76a9-490a-a5ce-5ba95dc756e3'
$ws.Range("C19").Value2 = '2024-07-03T09:52:00.000Z'
$ws.Range("D19").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I19").Value2 = 'https://www.notion.so/18-d37aca6076a9490aa5ce5ba95dc756e3'

# Row 20: 4 cell(s) changed
$ws.Range("B20").Value2 = 'e0bee04d-4ecb-4601-928c-0d8895b609ff'
$ws.Range("C20").Value2 = '2024-07-02T02:58:00.000Z'
$ws.Range("D20").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I20").Value2 = 'https://www.notion.so/19-e0bee04d4ecb4601928c0d8895b609ff'

# Row 21: 4 cell(s) changed
$ws.Range("B21").Value2 = '06f84bee-3de0-4d8e-a89d-0d662e49da6f'
$ws.Range("C21").Value2 = '2024-07-01T07:45:00.000Z'
$ws.Range("D21").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I21").Value2 = 'https://www.notion.so/20-06f84bee3de04d8ea89d0d662e49da6f'

# Row 22: 4 cell(s) changed
$ws.Range("B22").Value2 = '9c68fce3-7d4f-4432-a9b2-e94a31b35cb8'
$ws.Range("C22").Value2 = '2024-06-30T09:08:00.000Z'
$ws.Range("D22").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I22").Value2 = 'https://www.notion.so/21-9c68fce37d4f4432a9b2e94a31b35cb8'

# Row 23: 4 cell(s) changed
$ws.Range("B23").Value2 = 'f4840e29-9c1d-43c3-837e-0d16a331afc5'
$ws.Range("C23").Value2 = '2024-06-29T09:34:00.000Z'
$ws.Range("D23").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I23").Value2 = 'https://www.notion.so/22-f4840e299c1d43c3837e0d16a331afc5'

# Row 24: 4 cell(s) changed
$ws.Range("B24").Value2 = '413abd2b-d960-425d-830b-670598ce73e3'
$ws.Range("C24").Value2 = '2024-06-27T05:20:00.000Z'
$ws.Range("D24").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I24").Value2 = 'https://www.notion.so/23-413abd2bd960425d830b670598ce73e3'

# Row 25: 4 cell(s) changed
$ws.Range("B25").Value2 = '7f65d6ca-1a41-47e3-90c3-1c67b654a580'
$ws.Range("C25").Value2 = '2024-06-26T11:55:00.000Z'
$ws.Range("D25").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I25").Value2 = 'https://www.notion.so/24-7f65d6ca1a4147e390c31c67b654a580'

# Row 26: 3 cell(s) changed
$ws.Range("B26").Value2 = '5368b5a8-2d91-4fc3-a5f1-94c7162f5048'
$ws.Range("D26").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I26").Value2 = 'https://www.notion.so/25-5368b5a82d914fc3a5f194c7162f5048'

# Row 27: 3 cell(s) changed
$ws.Range("B27").Value2 = '8701f2d2-93c3-4e00-8567-2e3e7631cdb8'
$ws.Range("D27").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I27").Value2 = 'https://www.notion.so/26-8701f2d293c34e0085672e3e7631cdb8'

# Row 28: 3 cell(s) changed
$ws.Range("B28").Value2 = '456dae27-3fe8-4d68-bdf8-4f10f987d889'
$ws.Range("D28").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I28").Value2 = 'https://www.notion.so/27-456dae273fe84d68bdf84f10f987d889'

# Row 29: 3 cell(s) changed
$ws.Range("B29").Value2 = 'ccfa1087-8ca8-496b-8e93-2819f304be25'
$ws.Range("D29").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I29").Value2 = 'https://www.notion.so/28-ccfa10878ca8496b8e932819f304be25'

# Row 30: 3 cell(s) changed
$ws.Range("B30").Value2 = '1c8d6284-3aae-418c-9208-9c391cbbf64d'
$ws.Range("D30").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I30").Value2 = 'https://www.notion.so/29-1c8d62843aae418c92089c391cbbf64d'

# Row 31: 3 cell(s) changed
$ws.Range("B31").Value2 = 'ac86b813-bcdf-4436-856f-f6d3737b416c'
$ws.Range("D31").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I31").Value2 = 'https://www.notion.so/30-ac86b813bcdf4436856ff6d3737b416c'

# Row 32: 3 cell(s) changed
$ws.Range("B32").Value2 = '8e3a71e5-f927-4e4f-b215-d54c0188294e'
$ws.Range("D32").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I32").Value2 = 'https://www.notion.so/31-8e3a71e5f9274e4fb215d54c0188294e'

# Row 33: 3 cell(s) changed
$ws.Range("B33").Value2 = 'bc3f46b5-6495-4ce2-90fb-3abd2ad49c48'
$ws.Range("D33").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I33").Value2 = 'https://www.notion.so/32-bc3f46b564954ce290fb3abd2ad49c48'

# Row 34: 3 cell(s) changed
$ws.Range("B34").Value2 = '2dcf2fed-b2f3-45f7-8737-4acf8616b370'
$ws.Range("D34").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I34").Value2 = 'https://www.notion.so/33-2dcf2fedb2f345f787374acf8616b370'

# Row 35: 3 cell(s) changed
$ws.Range("B35").Value2 = '887b2441-0fb9-4d5d-86f3-59c23164bb5d'
$ws.Range("D35").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I35").Value2 = 'https://www.notion.so/34-887b24410fb94d5d86f359c23164bb5d'

# Row 36: 4 cell(s) changed
$ws.Range("B36").Value2 = 'b002fbe7-cff2-401f-b9ba-d07fd6cc561f'
$ws.Range("C36").Value2 = '2024-06-25T15:40:00.000Z'
$ws.Range("D36").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I36").Value2 = 'https://www.notion.so/35-b002fbe7cff2401fb9bad07fd6cc561f'

# Row 37: 3 cell(s) changed
$ws.Range("B37").Value2 = 'df6bda30-993c-4e7a-ab47-94d2c47b686c'
$ws.Range("D37").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I37").Value2 = 'https://www.notion.so/36-df6bda30993c4e7aab4794d2c47b686c'

# Row 38: 3 cell(s) changed
$ws.Range("B38").Value2 = '05ebb8d4-9dc1-4159-9836-122e738bf488'
$ws.Range("D38").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I38").Value2 = 'https://www.notion.so/37-05ebb8d49dc141599836122e738bf488'

# Row 39: 3 cell(s) changed
$ws.Range("B39").Value2 = '2ac05bc7-172c-4aa4-b654-799d67afe9e3'
$ws.Range("D39").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I39").Value2 = 'https://www.notion.so/38-2ac05bc7172c4aa4b654799d67afe9e3'

# Row 40: 3 cell(s) changed
$ws.Range("B40").Value2 = 'ffe6c779-91f7-4c50-b15d-d0453a138e44'
$ws.Range("D40").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I40").Value2 = 'https://www.notion.so/39-ffe6c77991f74c50b15dd0453a138e44'

# Row 41: 3 cell(s) changed
$ws.Range("B41").Value2 = '2b0f762d-d9a1-460a-b13f-2736f1b47f93'
$ws.Range("D41").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I41").Value2 = 'https://www.notion.so/40-2b0f762dd9a1460ab13f2736f1b47f93'

# Row 42: 3 cell(s) changed
$ws.Range("B42").Value2 = '8716b812-f61b-4f9c-90a0-142d28bdca52'
$ws.Range("D42").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I42").Value2 = 'https://www.notion.so/41-8716b812f61b4f9c90a0142d28bdca52'

# Row 43: 3 cell(s) changed
$ws.Range("B43").Value2 = '670caeaf-380a-49d7-80a2-f4774180a454'
$ws.Range("D43").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I43").Value2 = 'https://www.notion.so/42-670caeaf380a49d780a2f4774180a454'

# Row 44: 3 cell(s) changed
$ws.Range("B44").Value2 = '4159c13e-192d-4b3c-a97a-4b8514d8725c'
$ws.Range("D44").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I44").Value2 = 'https://www.notion.so/43-4159c13e192d4b3ca97a4b8514d8725c'

# Row 45: 3 cell(s) changed
$ws.Range("B45").Value2 = 'cf3e7eac-9355-4656-9394-e931803f630d'
$ws.Range("D45").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I45").Value2 = 'https://www.notion.so/44-cf3e7eac935546569394e931803f630d'

# Row 46: 3 cell(s) changed
$ws.Range("B46").Value2 = '9df73534-9716-4a8a-b249-39853ee869da'
$ws.Range("D46").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I46").Value2 = 'https://www.notion.so/45-9df7353497164a8ab24939853ee869da'

# Row 47: 3 cell(s) changed
$ws.Range("B47").Value2 = '5346b509-87b4-46f7-9918-eecb88e8cf6a'
$ws.Range("D47").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I47").Value2 = 'https://www.notion.so/46-5346b50987b446f79918eecb88e8cf6a'

# Row 48: 3 cell(s) changed
$ws.Range("B48").Value2 = 'a5122c43-3023-4baa-be9f-0dd89011734c'
$ws.Range("D48").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I48").Value2 = 'https://www.notion.so/47-a5122c4330234baabe9f0dd89011734c'

# Row 49: 3 cell(s) changed
$ws.Range("B49").Value2 = 'a532e7cf-9091-4f2e-9e8d-88f03ae7f109'
$ws.Range("D49").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I49").Value2 = 'https://www.notion.so/48-a532e7cf90914f2e9e8d88f03ae7f109'

# Row 50: 3 cell(s) changed
$ws.Range("B50").Value2 = 'ecbeca94-5cb6-43fd-b3ea-aa736eb36655'
$ws.Range("D50").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I50").Value2 = 'https://www.notion.so/49-ecbeca945cb643fdb3eaaa736eb36655'

# Row 51: 3 cell(s) changed
$ws.Range("B51").Value2 = '0cb9eb84-b38a-4c77-b4fb-c39974d76a86'
$ws.Range("D51").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I51").Value2 = 'https://www.notion.so/50-0cb9eb84b38a4c77b4fbc39974d76a86'

# Row 52: 3 cell(s) changed
$ws.Range("B52").Value2 = '71579228-ff4f-4c82-ab8d-e40170d28c84'
$ws.Range("D52").Value2 = '2024-07-25T15:02:00.000Z'
$ws.Range("I52").Value2 = 'https://www.notion.so/51-71579228ff4f4c82ab8de40170d28c84'

# Row 53: 3 cell(s) changed
$ws.Range("B53").Value2 = '808bf9c3-006a-4002-93f0-8edfcc7f0f77'
$ws.Range("D53").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I53").Value2 = 'https://www.notion.so/52-808bf9c3006a400293f08edfcc7f0f77'

# Row 54: 3 cell(s) changed
$ws.Range("B54").Value2 = 'f6aa15d5-81a7-44cd-b711-a7c11cd8091b'
$ws.Range("D54").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I54").Value2 = 'https://www.notion.so/53-f6aa15d581a744cdb711a7c11cd8091b'

# Row 55: 3 cell(s) changed
$ws.Range("B55").Value2 = '8f93fb5b-6434-414d-b158-5a6d8373e0cb'
$ws.Range("D55").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I55").Value2 = 'https://www.notion.so/54-8f93fb5b6434414db1585a6d8373e0cb'

# Row 56: 3 cell(s) changed
$ws.Range("B56").Value2 = '5f697f46-bc45-47b0-af9f-b79d2f89615c'
$ws.Range("D56").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I56").Value2 = 'https://www.notion.so/55-5f697f46bc4547b0af9fb79d2f89615c'

# Row 57: 3 cell(s) changed
$ws.Range("B57").Value2 = '0d51e83d-daba-4f7f-98a5-265e6f893acc'
$ws.Range("D57").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I57").Value2 = 'https://www.notion.so/56-0d51e83ddaba4f7f98a5265e6f893acc'

# Row 58: 3 cell(s) changed
$ws.Range("B58").Value2 = 'b6442ee6-2444-4519-ba61-ee66f711f543'
$ws.Range("D58").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I58").Value2 = 'https://www.notion.so/57-b6442ee624444519ba61ee66f711f543'

# Row 59: 3 cell(s) changed
$ws.Range("B59").Value2 = '5b30e0cb-e4f8-42e8-ac2f-3ca7b7a6aaaa'
$ws.Range("D59").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I59").Value2 = 'https://www.notion.so/58-5b30e0cbe4f842e8ac2f3ca7b7a6aaaa'

# Row 60: 3 cell(s) changed
$ws.Range("B60").Value2 = '7be7f5b9-f6a1-4c0f-91ec-5883dd1adf3b'
$ws.Range("D60").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I60").Value2 = 'https://www.notion.so/59-7be7f5b9f6a14c0f91ec5883dd1adf3b'

# Row 61: 3 cell(s) changed
$ws.Range("B61").Value2 = '5ae18489-d185-40db-bfc4-1263a08c863e'
$ws.Range("D61").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I61").Value2 = 'https://www.notion.so/60-5ae18489d18540dbbfc41263a08c863e'

# Row 62: 3 cell(s) changed
$ws.Range("B62").Value2 = '2cc2a873-7fd9-4612-a974-d6e1f3cfcad0'
$ws.Range("D62").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I62").Value2 = 'https://www.notion.so/61-2cc2a8737fd94612a974d6e1f3cfcad0'

# Row 63: 3 cell(s) changed
$ws.Range("B63").Value2 = '1bbc9ed6-97b3-463a-a707-3252641dd198'
$ws.Range("D63").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I63").Value2 = 'https://www.notion.so/62-1bbc9ed697b3463aa7073252641dd198'

# Row 64: 3 cell(s) changed
$ws.Range("B64").Value2 = 'c0b0d935-ef9a-45f7-879b-558e9ce4fcf3'
$ws.Range("D64").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I64").Value2 = 'https://www.notion.so/63-c0b0d935ef9a45f7879b558e9ce4fcf3'

# Row 65: 3 cell(s) changed
$ws.Range("B65").Value2 = 'de2a36cf-8771-44f0-a9a3-cd55abfe271b'
$ws.Range("D65").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I65").Value2 = 'https://www.notion.so/64-de2a36cf877144f0a9a3cd55abfe271b'

# Row 66: 3 cell(s) changed
$ws.Range("B66").Value2 = '52391eb5-b8b0-4d4c-a021-f5c63b61d834'
$ws.Range("D66").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I66").Value2 = 'https://www.notion.so/65-52391eb5b8b04d4ca021f5c63b61d834'

# Row 67: 3 cell(s) changed
$ws.Range("B67").Value2 = '68bb50dc-b720-47b7-9983-76df9b3742f6'
$ws.Range("D67").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I67").Value2 = 'https://www.notion.so/66-68bb50dcb72047b7998376df9b3742f6'

# Row 68: 3 cell(s) changed
$ws.Range("B68").Value2 = '58b1fa10-c5cf-4bcd-be34-669b0c698923'
$ws.Range("D68").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I68").Value2 = 'https://www.notion.so/67-58b1fa10c5cf4bcdbe34669b0c698923'

# Row 69: 3 cell(s) changed
$ws.Range("B69").Value2 = 'c3ede5dd-eeba-4a26-8753-8586469f6e8a'
$ws.Range("D69").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I69").Value2 = 'https://www.notion.so/68-c3ede5ddeeba4a2687538586469f6e8a'

# Row 70: 3 cell(s) changed
$ws.Range("B70").Value2 = '46806830-e1bd-4cb5-94a5-26deedae2a38'
$ws.Range("D70").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I70").Value2 = 'https://www.notion.so/69-46806830e1bd4cb594a526deedae2a38'

# Row 71: 3 cell(s) changed
$ws.Range("B71").Value2 = '2e1d33d8-b833-4cca-acee-1f03482d66eb'
$ws.Range("D71").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I71").Value2 = 'https://www.notion.so/70-2e1d33d8b8334ccaacee1f03482d66eb'

# Row 72: 3 cell(s) changed
$ws.Range("B72").Value2 = '96324ff3-7f59-48f0-b916-fc6536c5dade'
$ws.Range("D72").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I72").Value2 = 'https://www.notion.so/71-96324ff37f5948f0b916fc6536c5dade'

# Row 73: 3 cell(s) changed
$ws.Range("B73").Value2 = '0e7dfb28-8771-4894-8e46-4365dc6f2ad2'
$ws.Range("D73").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I73").Value2 = 'https://www.notion.so/72-0e7dfb28877148948e464365dc6f2ad2'

# Row 74: 3 cell(s) changed
$ws.Range("B74").Value2 = '843f0e90-37ff-4591-9e45-0424e4cbac9e'
$ws.Range("D74").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I74").Value2 = 'https://www.notion.so/73-843f0e9037ff45919e450424e4cbac9e'

# Row 75: 3 cell(s) changed
$ws.Range("B75").Value2 = '6049f3d7-42d0-4838-8a7a-1856688cc253'
$ws.Range("D75").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I75").Value2 = 'https://www.notion.so/74-6049f3d742d048388a7a1856688cc253'

# Row 76: 3 cell(s) changed
$ws.Range("B76").Value2 = '612cdb59-27f9-493b-b674-72aff3b5927b'
$ws.Range("D76").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I76").Value2 = 'https://www.notion.so/75-612cdb5927f9493bb67472aff3b5927b'

# Row 77: 3 cell(s) changed
$ws.Range("B77").Value2 = '27085ca3-d514-4473-a3bb-da984c52e7e4'
$ws.Range("D77").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I77").Value2 = 'https://www.notion.so/76-27085ca3d5144473a3bbda984c52e7e4'

# Row 78: 3 cell(s) changed
$ws.Range("B78").Value2 = '5be9f11c-c64f-4f61-bb71-e7d2fbf96855'
$ws.Range("D78").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I78").Value2 = 'https://www.notion.so/77-5be9f11cc64f4f61bb71e7d2fbf96855'

# Row 79: 3 cell(s) changed
$ws.Range("B79").Value2 = '3f9e3e71-188e-48d5-997a-71070175e42b'
$ws.Range("D79").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I79").Value2 = 'https://www.notion.so/78-3f9e3e71188e48d5997a71070175e42b'

# Row 80: 4 cell(s) changed
$ws.Range("B80").Value2 = 'c4a99c5c-f44e-4039-a045-8ad12b4f6046'
$ws.Range("C80").Value2 = '2024-06-25T15:39:00.000Z'
$ws.Range("D80").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I80").Value2 = 'https://www.notion.so/79-c4a99c5cf44e4039a0458ad12b4f6046'

# Row 81: 3 cell(s) changed
$ws.Range("B81").Value2 = '214991f6-4649-4c4d-b595-7aa6576c69ad'
$ws.Range("D81").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I81").Value2 = 'https://www.notion.so/80-214991f646494c4db5957aa6576c69ad'

# Row 82: 3 cell(s) changed
$ws.Range("B82").Value2 = 'd9f9ae69-896d-4f76-b572-08645f10e3fa'
$ws.Range("D82").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I82").Value2 = 'https://www.notion.so/81-d9f9ae69896d4f76b57208645f10e3fa'

# Row 83: 3 cell(s) changed
$ws.Range("B83").Value2 = '73f7f862-f114-41db-9221-0fde8f9c580d'
$ws.Range("D83").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I83").Value2 = 'https://www.notion.so/82-73f7f862f11441db92210fde8f9c580d'

# Row 84: 3 cell(s) changed
$ws.Range("B84").Value2 = '631d954b-83a8-47d4-8683-983aad0b002d'
$ws.Range("D84").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I84").Value2 = 'https://www.notion.so/83-631d954b83a847d48683983aad0b002d'

# Row 85: 5 cell(s) changed
$ws.Range("B85").Value2 = '3e313c45-1109-4d19-96f7-3d46baff66af'
$ws.Range("C85").Value2 = '2024-06-25T15:36:00.000Z'
$ws.Range("D85").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I85").Value2 = 'https://www.notion.so/84-3e313c4511094d1996f73d46baff66af'
$ws.Range("L85").Value2 = '41cabcaf-915d-46a5-8eff-38727be27269'

# Row 86: 3 cell(s) changed
$ws.Range("B86").Value2 = '15dc290b-de0b-4ca9-ba77-0946c5f13582'
$ws.Range("D86").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I86").Value2 = 'https://www.notion.so/85-15dc290bde0b4ca9ba770946c5f13582'

# Row 87: 3 cell(s) changed
$ws.Range("B87").Value2 = '85a1185c-5d12-479f-808c-48ad22a864c1'
$ws.Range("D87").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("I87").Value2 = 'https://www.notion.so/86-85a1185c5d12479f808c48ad22a864c1'

# Row 88: 47 cell(s) changed
$ws.Range("A88").Value2 = 'page'
$ws.Range("B88").Value2 = 'b4c82aa7-7e7f-4289-b0ba-1b14e683031d'
$ws.Range("C88").Value2 = '2024-06-25T08:24:00.000Z'
$ws.Range("D88").Value2 = '2024-07-25T15:03:00.000Z'
$ws.Range("E88").Value2 = ''
$ws.Range("F88").Value2 = ''
$ws.Range("G88").Value2 = $false
$ws.Range("H88").Value2 = $false
$ws.Range("I88").Value2 = 'https://www.notion.so/87-b4c82aa77e7f4289b0ba1b14e683031d'
$ws.Range("J88").Value2 = ''
$ws.Range("K88").Value2 = 'user'
$ws.Range("L88").Value2 = '532a166e-c2d9-42ff-bed3-a363f43543fb'
$ws.Range("M88").Value2 = 'user'
$ws.Range("N88").Value2 = '41cabcaf-915d-46a5-8eff-38727be27269'
$ws.Range("O88").Value2 = 'database_id'
$ws.Range("P88").Value2 = '95908de9-fba9-4247-8bcd-04b57e56bd1b'
$ws.Range("Q88").Value2 = 'JgG%7B'
$ws.Range("R88").Value2 = 'number'
$ws.Range("S88").Value2 = 0
$ws.Range("T88").Value2 = 'Nk%3CD'
$ws.Range("U88").Value2 = 'number'
$ws.Range("V88").Value2 = 0
$ws.Range("W88").Value2 = 'Q%5B%5Bp'
$ws.Range("X88").Value2 = 'formula'
$ws.Range("Y88").Value2 = 'number'
$ws.Range("Z88").Value2 = 0
$ws.Range("AA88").Value2 = 'Vds%5C'
$ws.Range("AB88").Value2 = 'number'
$ws.Range("AC88").Value2 = 2050000
$ws.Range("AD88").Value2 = 'bwGG'
$ws.Range("AE88").Value2 = 'number'
$ws.Range("AF88").Value2 = 0
$ws.Range("AG88").Value2 = 'kciy'
$ws.Range("AH88").Value2 = 'date'
$ws.Range("AI88").Value2 = '2024-07-25'
$ws.Range("AJ88").Value2 = ''
$ws.Range("AK88").Value2 = ''
$ws.Range("AL88").Value2 = 'nWGB'
$ws.Range("AM88").Value2 = 'formula'
$ws.Range("AN88").Value2 = 'number'
$ws.Range("AO88").Value2 = -2050000
$ws.Range("AP88").Value2 = 'pq%3Ci'
$ws.Range("AQ88").Value2 = 'number'
$ws.Range("AR88").Value2 = 0
$ws.Range("AS88").Value2 = 'title'
$ws.Range("AT88").Value2 = 'title'
$ws.Range("AU88").Value2 = '[{''type'': ''text'', ''text'': {''content'': ''87'', ''link'': None}, ''annotations'': {''bold'': False, ''italic'': False, ''strikethrough'': False, ''underline'': False, ''code'': False, ''color'': ''default''}, ''plain_text'': ''87'', ''href'': None}]'
